$d = $word.ActiveDocument

# "Versi" + "on" (two separate runs) -> merge into a single "Version" run.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# " 2" -> " 1." (the trailing period moves here from the run that follows
# the _GoBack bookmark).
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# Remove the now-redundant standalone "." run that used to sit after the
# _GoBack bookmark.
$fullText = $d.Content.Text
$dotIndex = $fullText.Length - 2
$r = $d.Range($dotIndex, $dotIndex + 1)
if ($r.Text -eq ".") {
    $r.Delete()
}
